# Auto-generated edit script: updates currentAveragePrice / HQ / profit columns
# (H..N) for specific leve rows across multiple crafting-job sheets, per scheduled
# Moogle Profits price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 40
$ws.Range("H40").Value = 5040.353
$ws.Range("J40").Value = 7613.7144
$ws.Range("L40").Value = 7613.7144
$ws.Range("N40").Value = -7963.7144

# ALC row 113
$ws.Range("H113").Value = 5622
$ws.Range("I113").Value = 4123
$ws.Range("K113").Value = 4123
$ws.Range("M113").Value = -869

# ALC row 132
$ws.Range("H132").Value = 2817.6758
$ws.Range("I132").Value = 2170.9412
$ws.Range("J132").Value = 10147.333
$ws.Range("K132").Value = 6512.823600000001
$ws.Range("L132").Value = 30441.999
$ws.Range("M132").Value = -3982.823600000001
$ws.Range("N132").Value = -35501.999

# ALC row 137
$ws.Range("H137").Value = 2583.96
$ws.Range("I137").Value = 2399.9546
$ws.Range("K137").Value = 7199.8638
$ws.Range("M137").Value = -4649.8638

# ALC row 138
$ws.Range("H138").Value = 2048.1133
$ws.Range("J138").Value = 3761.9375
$ws.Range("L138").Value = 11285.8125
$ws.Range("N138").Value = -21565.8125


$ws = $wb.Worksheets.Item("BSM")
# BSM row 64
$ws.Range("H64").Value = 750
$ws.Range("I64").Value = 750
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 750
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -525
$ws.Range("N64").ClearContents()

# BSM row 67
$ws.Range("H67").Value = 750
$ws.Range("I67").Value = 750
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 750
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 30
$ws.Range("N67").ClearContents()

# BSM row 94
$ws.Range("H94").Value = 540.25
$ws.Range("I94").Value = 422.2
$ws.Range("K94").Value = 422.2
$ws.Range("M94").Value = 28.80000000000001


$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 16721.908
$ws.Range("I31").Value = 3000
$ws.Range("J31").Value = 18094.1
$ws.Range("K31").Value = 3000
$ws.Range("L31").Value = 18094.1
$ws.Range("M31").Value = -2705
$ws.Range("N31").Value = -18684.1

# CRP row 34
$ws.Range("H34").Value = 16721.908
$ws.Range("I34").Value = 3000
$ws.Range("J34").Value = 18094.1
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 18094.1
$ws.Range("M34").Value = -2798
$ws.Range("N34").Value = -18498.1

# CRP row 50
$ws.Range("H50").Value = 51096.11
$ws.Range("I50").Value = 43995
$ws.Range("J50").Value = 59972.5
$ws.Range("K50").Value = 43995
$ws.Range("L50").Value = 59972.5
$ws.Range("M50").Value = -43370
$ws.Range("N50").Value = -61222.5

# CRP row 51
$ws.Range("H51").Value = 39949.25
$ws.Range("J51").Value = 49932.668
$ws.Range("L51").Value = 49932.668
$ws.Range("N51").Value = -51404.668

# CRP row 60
$ws.Range("H60").Value = 31586
$ws.Range("J60").Value = 46275.5
$ws.Range("L60").Value = 46275.5
$ws.Range("N60").Value = -47297.5

# CRP row 61
$ws.Range("H61").Value = 39949.25
$ws.Range("J61").Value = 49932.668
$ws.Range("L61").Value = 49932.668
$ws.Range("N61").Value = -50628.668

# CRP row 132
$ws.Range("H132").Value = 3386.3215
$ws.Range("I132").Value = 2455.4546
$ws.Range("K132").Value = 7366.3638
$ws.Range("M132").Value = -4836.3638


$ws = $wb.Worksheets.Item("CUL")
# CUL row 12
$ws.Range("H12").Value = 355.54544
$ws.Range("I12").Value = 103
$ws.Range("J12").Value = 450.25
$ws.Range("K12").Value = 309
$ws.Range("L12").Value = 1350.75
$ws.Range("M12").Value = -136
$ws.Range("N12").Value = -1696.75

# CUL row 38
$ws.Range("H38").Value = 93.7
$ws.Range("I38").Value = 94.333336
$ws.Range("J38").Value = 92.75
$ws.Range("K38").Value = 283.000008
$ws.Range("L38").Value = 278.25
$ws.Range("M38").Value = 63.99999200000002
$ws.Range("N38").Value = -972.25

# CUL row 68
$ws.Range("H68").Value = 2174
$ws.Range("J68").Value = 2666.6667
$ws.Range("L68").Value = 8000.000100000001
$ws.Range("N68").Value = -9622.000100000001

# CUL row 71
$ws.Range("H71").Value = 2174
$ws.Range("J71").Value = 2666.6667
$ws.Range("L71").Value = 24000.0003
$ws.Range("N71").Value = -32112.0003

# CUL row 92
$ws.Range("H92").Value = 4599.5
$ws.Range("J92").Value = 4949
$ws.Range("L92").Value = 14847
$ws.Range("N92").Value = -17343

# CUL row 141
$ws.Range("H141").Value = 6214
$ws.Range("I141").Value = 5599.75
$ws.Range("J141").Value = 7033
$ws.Range("K141").Value = 16799.25
$ws.Range("L141").Value = 21099
$ws.Range("M141").Value = -11619.25
$ws.Range("N141").Value = -31459


$ws = $wb.Worksheets.Item("GSM")
# GSM row 97
$ws.Range("H97").Value = 420.875
$ws.Range("I97").Value = 366.56522
$ws.Range("K97").Value = 366.56522
$ws.Range("M97").Value = 129.43478

# GSM row 132
$ws.Range("H132").Value = 7433.227
$ws.Range("I132").Value = 6676.55
$ws.Range("K132").Value = 20029.65
$ws.Range("M132").Value = -17499.65


$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 9499.462
$ws.Range("I16").Value = 2099.3
$ws.Range("K16").Value = 2099.3
$ws.Range("M16").Value = -1929.3

# LTW row 40
$ws.Range("H40").Value = 6638.9
$ws.Range("I40").Value = 5965
$ws.Range("J40").Value = 7649.75
$ws.Range("K40").Value = 5965
$ws.Range("L40").Value = 7649.75
$ws.Range("M40").Value = -5829
$ws.Range("N40").Value = -7921.75

# LTW row 61
$ws.Range("H61").Value = 1700.7778
$ws.Range("I61").Value = 729.8570999999999
$ws.Range("J61").Value = 5099
$ws.Range("K61").Value = 729.8570999999999
$ws.Range("L61").Value = 5099
$ws.Range("M61").Value = -527.8570999999999
$ws.Range("N61").Value = -5503

# LTW row 93
$ws.Range("H93").Value = 1296.4667
$ws.Range("I93").Value = 1001.7273
$ws.Range("J93").Value = 2107
$ws.Range("K93").Value = 1001.7273
$ws.Range("L93").Value = 2107
$ws.Range("M93").Value = 246.2727
$ws.Range("N93").Value = -4603

# LTW row 113
$ws.Range("H113").Value = 1700.7778
$ws.Range("I113").Value = 729.8570999999999
$ws.Range("J113").Value = 5099
$ws.Range("K113").Value = 729.8570999999999
$ws.Range("L113").Value = 5099
$ws.Range("M113").Value = 1440.1429
$ws.Range("N113").Value = -9439

# LTW row 122
$ws.Range("H122").Value = 5738.5835
$ws.Range("I122").Value = 4345.6665
$ws.Range("J122").Value = 6202.8887
$ws.Range("K122").Value = 13036.9995
$ws.Range("L122").Value = 18608.6661
$ws.Range("M122").Value = -10586.9995
$ws.Range("N122").Value = -23508.6661


$ws = $wb.Worksheets.Item("WVR")
# WVR row 32
$ws.Range("H32").Value = 11805.2
$ws.Range("I32").Value = 12256.5
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 12256.5
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -11939.5
$ws.Range("N32").Value = -10634

# WVR row 38
$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 10000
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -9527
$ws.Range("N38").ClearContents()

# WVR row 122
$ws.Range("H122").Value = 3789
$ws.Range("I122").Value = 3918.842
$ws.Range("K122").Value = 11756.526
$ws.Range("M122").Value = -9306.526

